# [BI-1613] Update TAF to include term type
# Add a new "Term Type" column (R) to the trait-import template header row,
# matching the style already used for the other header cells (bold,
# wrap-text, bordered).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

$header = $ws.Range("R1")
$header.Value = "Term Type"
$header.Font.Bold = $true
$header.Font.Size = 11
$header.WrapText = $true
$header.Borders.LineStyle = 1

$header.Select() | Out-Null
